# Consolidate the trailing "Results" column into column I (was J) and drop
# the now-empty spacer column. Mirrors the manual edit of right-clicking
# column J -> Delete (which shifts J's content into I and removes the
# column entirely), then re-sizing column H to fit its header text and
# leaving the freshly-vacated column I selected.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Move J1's ("Results") content into I1 (which was just an empty, filled
# placeholder cell) before the column shift swallows it.
$ws.Range("I1").Value2 = $ws.Range("J1").Value2

# Delete the now-redundant column J entirely; everything to its right (none,
# here) shifts left and the used range shrinks from A1:J9 to A1:I9.
$ws.Range("J1").EntireColumn.Delete()

# Column H ("Sample.Unit.Modifier" header) needs a fitted custom width of
# 18 characters. Feeding 17.1666... through ColumnWidth lands exactly on a
# stored width of 18 once Excel's char->pixel->width rounding is applied.
$ws.Columns("H").ColumnWidth = 17.166666666666668

# Leave the whole of column I selected, matching the post-edit selection.
$ws.Range("I1").EntireColumn.Select()
